$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 295, pushing the existing rows 295-320 down to 296-321.
$ws.Rows.Item(295).Insert()

# Populate the newly-inserted row 295 with the new weekly record.
$ws.Cells.Item(295, 1).Value = 5
$ws.Cells.Item(295, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(295, 3).Value = "Maule"
$ws.Cells.Item(295, 4).Value = 44578
$ws.Cells.Item(295, 5).Value = 7
$ws.Cells.Item(295, 6).Value = "Fruta"
$ws.Cells.Item(295, 7).Value = 100109
$ws.Cells.Item(295, 8).Value = "Uva"
$ws.Cells.Item(295, 9).Value = 100109001
$ws.Cells.Item(295, 10).Value = "Uva"
$ws.Cells.Item(295, 11).Value = "Superior Seedless"
$ws.Cells.Item(295, 12).Value = "Primera"
$ws.Cells.Item(295, 13).Value = 300
$ws.Cells.Item(295, 14).Value = 10000
$ws.Cells.Item(295, 15).Value = 10000
$ws.Cells.Item(295, 16).Value = 10000
$ws.Cells.Item(295, 17).Value = "$/bandeja 9 kilos"
$ws.Cells.Item(295, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(295, 19).Value = 1111
$ws.Cells.Item(295, 20).Value = 9
